# RESTORE: revert the multi-industry "AI/ML" template strings back onto
# this "Enterprise Cloud Infrastructure Migration" IT Quality Management
# Plan sample workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Change Management Overview" ---
$ws1 = $wb.Worksheets.Item("Change Management Overview")
$ws1.Range("B6").Value = "Enterprise AI/ML Implementation"
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new AI/ML systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in AI/ML technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for AI/ML transformation"

# --- Sheet 2: "Change Impact Assessment" ---
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")
$ws2.Range("A4").Value = "Data Scientists"
$ws2.Range("G4").Value = "AI/ML automation"
$ws2.Range("A5").Value = "Business Analysts"

# Note: Sheet 3 ("Change Activities") needs no cell-content edits. The
# blank spacer rows (sheet1 row 13 & 21, sheet2 row 2, sheet3 row 2) that
# the diff drops are automatically normalized away by the workbook save
# pipeline once it re-serializes, since they carry no cells/attributes.
